$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.938.58'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '3.397.18'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '579.42'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('D6').Value = '177.11'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').Value = '0.196'
$ws.Range('E9').Value = '  +6.77%  '
$ws.Range('D10').Value = '0.582'
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('D11').Value = '47.96'
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '0.0000280'
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('D13').Value = '675.93'
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('D14').Value = '3.947.78'
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('D15').Value = '8.57'
$ws.Range('E15').Value = '  +1.30%  '
$ws.Range('D16').Value = '69.122.28'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '3.398.02'
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '17.70'
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').Value = '11.24'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').Value = '0.906'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').Value = '5.35'
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('D23').Value = '16.94'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '100.41'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = '3.88'
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('D26').Value = '2.67'
$ws.Range('E26').Value = '  -1.01%  '
$ws.Range('D27').Value = '9.61'
$ws.Range('E27').Value = '  +1.28%  '
$ws.Range('D28').Value = '33.33'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('D29').Value = '8.69'
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('D30').Value = '6.82'
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('D31').Value = '3.69'
$ws.Range('E31').Value = '  +7.84%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '10.95'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '549.58'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '3.597.18'
$ws.Range('E37').Value = '  -3.31%  '
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('D39').Value = '34.82'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').Value = '0.0₃0731'
$ws.Range('E40').Value = '  +8.97%  '
$ws.Range('D41').Value = '3.26'
$ws.Range('E41').Value = '  +2.75%  '
$ws.Range('D42').Value = '2.66'
$ws.Range('E42').Value = '  +1.95%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').Value = '3.36'
$ws.Range('E43').Value = '  +3.05%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0423'
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.332'
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = '2.65'
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.128'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '1.39'
$ws.Range('E48').Value = '  +3.51%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '130.84'
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').Value = '2.68'
$ws.Range('E51').Value = '  +3.85%  '
